$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# Copy the header style (bold, centered, bordered) from H1 onto the new headers.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

for ($r = 2; $r -le 27; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $hVal
}
